$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Global")

# Copy formatting (styles/borders) of the last data row (row 4) down into the
# new row 5 that records the out-of-stock "XL" item.
$ws.Range("A4:G4").Copy($ws.Range("A5:G5"))

# Row 5 data: BusinessProcess=3, Product="Ray Packable", Size="XL" (out of stock).
$ws.Cells.Item(5, 5).Value2 = 3
$ws.Cells.Item(5, 6).Value2 = "Ray Packable"
$ws.Cells.Item(5, 7).Value2 = "XL"

Write-Output "done"
